# Correcciones en reglas del documento stock actual
# Applies the "Semana 06" stock-rule corrections: several SKUs get their
# computed order/stock figures recalculated and are hidden from the visible
# report (their rule no longer triggers an order), a couple of swapped size
# rows (10/12 and 15/17) get their data swapped back to the correct SKU, and
# the summary metrics at the bottom of the sheet are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - 9203010004: no longer ordered this week, row hidden
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = 0
$ws.Range("U3").Value = 0
$ws.Rows.Item(3).Hidden = $true

# Row 8 - 9402010010: quantities halved
$ws.Range("M8").Value = 3.4
$ws.Range("N8").Value = 2.04
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("U8").Value = 1

# Row 9 - 9401010010: no longer ordered this week, row hidden
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("P9").Value = 8
$ws.Range("Q9").Value = 0
$ws.Range("U9").Value = 0
$ws.Rows.Item(9).Hidden = $true

# Row 10 - 9401010025: swapped with row 12 (now the 11CM/C variant), hidden
$ws.Range("C10").Value = "11CM     "
$ws.Range("G10").Value = 0.58
$ws.Range("H10").Value = 0.23
$ws.Range("I10").Value = "C"
$ws.Range("J10").Value = "REDUCIR 40%"
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("P10").Value = 40
$ws.Range("Q10").Value = 0
$ws.Range("S10").Value = 20
$ws.Range("T10").Value = 20
$ws.Range("U10").Value = 0
$ws.Rows.Item(10).Hidden = $true

# Row 11 - 9401010025: no longer ordered this week, row hidden
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("P11").Value = 24
$ws.Range("Q11").Value = 0
$ws.Range("U11").Value = 0
$ws.Rows.Item(11).Hidden = $true

# Row 12 - 9401010025: swapped with row 10 (now the 35CM/A variant), stays visible
$ws.Range("C12").Value = "35CM     "
$ws.Range("G12").Value = 5.98
$ws.Range("H12").Value = 2.39
$ws.Range("I12").Value = "A"
$ws.Range("J12").Value = "REDUCIR 19%"
$ws.Range("P12").Value = 18
$ws.Range("S12").Value = 0
$ws.Range("T12").Value = 0

# Row 15 - 9401020010: swapped with row 17 (now the 15CM variant), hidden
$ws.Range("C15").Value = "15CM     "
$ws.Range("G15").Value = 0.58
$ws.Range("H15").Value = 0.23
$ws.Range("J15").Value = "AUMENTAR 19%"
$ws.Range("L15").Value = -3
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 0
$ws.Range("U15").Value = 0
$ws.Rows.Item(15).Hidden = $true

# Row 16 - 9401020010: no longer ordered this week, row hidden
$ws.Range("L16").Value = -5
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 0
$ws.Range("U16").Value = 0
$ws.Rows.Item(16).Hidden = $true

# Row 17 - 9401020010: swapped with row 15 (now the 26CM variant), hidden
$ws.Range("C17").Value = "26CM     "
$ws.Range("G17").Value = 1.25
$ws.Range("H17").Value = 0.5
$ws.Range("J17").Value = "REDUCIR 19%"
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("P17").Value = 10
$ws.Range("Q17").Value = 0
$ws.Range("U17").Value = 0
$ws.Rows.Item(17).Hidden = $true

# Row 18 - 9401020012: quantities rescaled
$ws.Range("M18").Value = 9.449999999999999
$ws.Range("N18").Value = 5.67
$ws.Range("P18").Value = 7
$ws.Range("Q18").Value = 0
$ws.Range("U18").Value = 21

# Row 19 - 9401050006: no longer ordered this week, row hidden
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("P19").Value = 15
$ws.Range("Q19").Value = 0
$ws.Range("U19").Value = 0
$ws.Rows.Item(19).Hidden = $true

# Row 20 - 9401070015: no longer ordered this week, row hidden
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("P20").Value = 18
$ws.Range("Q20").Value = 0
$ws.Range("U20").Value = 0
$ws.Rows.Item(20).Hidden = $true

# Row 21 - 9201050008: no longer ordered this week, row hidden
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("P21").Value = 152
$ws.Range("Q21").Value = 0
$ws.Range("S21").Value = 6
$ws.Range("T21").Value = 6
$ws.Range("U21").Value = 0
$ws.Rows.Item(21).Hidden = $true

# Row 22 - 9201050009: no longer ordered this week, row hidden
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("P22").Value = 214
$ws.Range("Q22").Value = 0
$ws.Range("S22").Value = 9
$ws.Range("T22").Value = 9
$ws.Range("U22").Value = 0
$ws.Rows.Item(22).Hidden = $true

# Row 23 - 9201050013: no longer ordered this week, row hidden
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("P23").Value = 14
$ws.Range("Q23").Value = 0
$ws.Range("U23").Value = 0
$ws.Rows.Item(23).Hidden = $true

# Row 24 - 9701010001: no longer ordered this week, row hidden
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("P24").Value = 17
$ws.Range("Q24").Value = 0
$ws.Range("U24").Value = 0
$ws.Rows.Item(24).Hidden = $true

# Row 26 - 9402010044: diferencia stock corrected
$ws.Range("L26").Value = 1

# Row 27 - 9402020004: no longer ordered this week, row hidden
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("P27").Value = 57
$ws.Range("Q27").Value = 0
$ws.Range("U27").Value = 0
$ws.Rows.Item(27).Hidden = $true

# Row 29 - 9201050010: no longer ordered this week, row hidden
$ws.Range("L29").Value = -1
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("P29").Value = 78
$ws.Range("Q29").Value = 0
$ws.Range("S29").Value = 8
$ws.Range("T29").Value = 8
$ws.Range("U29").Value = 0
$ws.Rows.Item(29).Hidden = $true

# Row 31 - 9201050012: no longer ordered this week, row hidden
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("P31").Value = 130
$ws.Range("Q31").Value = 0
$ws.Range("U31").Value = 0
$ws.Rows.Item(31).Hidden = $true

# Row 34 - 9104020017: quantities halved
$ws.Range("M34").Value = 32.88
$ws.Range("N34").Value = 19.73
$ws.Range("P34").Value = 1
$ws.Range("Q34").Value = 1
$ws.Range("U34").Value = 1

# Row 35 - 9406010020: no longer ordered this week, row hidden
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("P35").Value = 4
$ws.Range("Q35").Value = 0
$ws.Range("S35").Value = 1
$ws.Range("T35").Value = 1
$ws.Range("U35").Value = 0
$ws.Rows.Item(35).Hidden = $true

# Row 36 - 9201040033: no longer ordered this week, row hidden
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("P36").Value = 118
$ws.Range("Q36").Value = 0
$ws.Range("S36").Value = 3
$ws.Range("T36").Value = 3
$ws.Range("U36").Value = 0
$ws.Rows.Item(36).Hidden = $true

# Row 37 - 9201040034: no longer ordered this week, row hidden
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("P37").Value = 60
$ws.Range("Q37").Value = 0
$ws.Range("U37").Value = 0
$ws.Rows.Item(37).Hidden = $true

# Row 38 - 9303080010: no longer ordered this week, row hidden
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("P38").Value = 2
$ws.Range("Q38").Value = 0
$ws.Range("U38").Value = 0
$ws.Rows.Item(38).Hidden = $true

# Summary metrics block (rows 40-52)
$ws.Range("C41").Value = 58

# C43 holds a currency-look-alike literal string ("1271.48€"). Assigning it
# directly would get auto-converted to a numeric currency value by Excel's
# smart-typing, so force literal text with a leading apostrophe and then
# restore the original (unmodified) cell formatting via a format-only paste
# from the untouched sibling cell C44, which carries the same base style.
$ws.Range("C43").Value = "'1271.48€"
$ws.Range("C44").Copy() | Out-Null
$ws.Range("C43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C52").Value = -6
